$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33, pushing existing rows 33:64 down to 34:65
# (weekly data refresh - a new observation is prepended to this block).
$ws.Rows(33).Insert()

# Populate the newly inserted row 33 with the latest weekly observation.
$ws.Cells.Item(33, 1).Value = 4
$ws.Cells.Item(33, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(33, 3).Value = "Los Lagos"
$ws.Cells.Item(33, 4).Value = 44638
$ws.Cells.Item(33, 5).Value = 10
$ws.Cells.Item(33, 6).Value = 100112031
$ws.Cells.Item(33, 7).Value = "Poroto verde"
$ws.Cells.Item(33, 8).Value = "Magnum"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 45
$ws.Cells.Item(33, 11).Value = 27000
$ws.Cells.Item(33, 12).Value = 27000
$ws.Cells.Item(33, 13).Value = 27000
$ws.Cells.Item(33, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(33, 15).Value = "Región Metropolitana"
$ws.Cells.Item(33, 16).Value = 1080
$ws.Cells.Item(33, 17).Value = 25
$ws.Cells.Item(33, 18).Value = "Hortaliza"
